# Generate Report for Handback
# Marks the two localized files as handed back (status + handback file/date)
# on the Overview sheet and on each language sheet (zh-cn, de-de).

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---- Overview sheet: update the Status-like columns (zh-cn / de-de) ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $newStatus
$overview.Range("C2").Value = $newStatus
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B2").Value = $newStatus
$zh.Range("B3").Value = $newStatus

# Row 2: 653de435-9040-4ca6-864d-6e5c29891627.md
$zh.Range("E2").Value = "653de435-9040-4ca6-864d-6e5c29891627.md"
$zh.Range("E2").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed9cb24aaeb09db6c3deebd985d5a2018c685bc3/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")

$zh.Range("F2").Value = "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf"
$zh.Range("F2").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63bcb5d53137b88e3d1a293772ea166b7e0f770b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/gt/653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf", "", "", "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.zh-cn.xlf")

$zh.Range("G2").Value = "2016-02-18 09:51:02"
$zh.Range("H2").Value = "Include"

# Row 3: a634b5f3-a252-4698-b996-c9ad1c439b66.md
$zh.Range("E3").Value = "a634b5f3-a252-4698-b996-c9ad1c439b66.md"
$zh.Range("E3").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed9cb24aaeb09db6c3deebd985d5a2018c685bc3/e2e/a634b5f3-a252-4698-b996-c9ad1c439b66.md", "", "", "a634b5f3-a252-4698-b996-c9ad1c439b66.md")

$zh.Range("F3").Value = "a634b5f3-a252-4698-b996-c9ad1c439b66.1c5614d47b2feca54d8e05ce9408c9d1c2230969.zh-cn.xlf"
$zh.Range("F3").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63bcb5d53137b88e3d1a293772ea166b7e0f770b/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/gt/a634b5f3-a252-4698-b996-c9ad1c439b66.1c5614d47b2feca54d8e05ce9408c9d1c2230969.zh-cn.xlf", "", "", "a634b5f3-a252-4698-b996-c9ad1c439b66.1c5614d47b2feca54d8e05ce9408c9d1c2230969.zh-cn.xlf")

$zh.Range("G3").Value = "2016-02-18 09:51:02"
$zh.Range("H3").Value = "Include"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("B2").Value = $newStatus
$de.Range("B3").Value = $newStatus

# Row 2: 653de435-9040-4ca6-864d-6e5c29891627.md
$de.Range("E2").Value = "653de435-9040-4ca6-864d-6e5c29891627.md"
$de.Range("E2").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/ed9cb24aaeb09db6c3deebd985d5a2018c685bc3/e2e/653de435-9040-4ca6-864d-6e5c29891627.md", "", "", "653de435-9040-4ca6-864d-6e5c29891627.md")

$de.Range("F2").Value = "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf"
$de.Range("F2").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c3cc36ee36b707a5742c80145a4eade4972fa2c9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/gt/653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf", "", "", "653de435-9040-4ca6-864d-6e5c29891627.dcb53bc46cef1baa3842aff3ceef28a8e9bffd2f.de-de.xlf")

$de.Range("G2").Value = "2016-02-18 09:51:25"

# Row 3: a634b5f3-a252-4698-b996-c9ad1c439b66.md
$de.Range("E3").Value = "a634b5f3-a252-4698-b996-c9ad1c439b66.md"
$de.Range("E3").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/ed9cb24aaeb09db6c3deebd985d5a2018c685bc3/e2e/a634b5f3-a252-4698-b996-c9ad1c439b66.md", "", "", "a634b5f3-a252-4698-b996-c9ad1c439b66.md")

$de.Range("F3").Value = "a634b5f3-a252-4698-b996-c9ad1c439b66.1c5614d47b2feca54d8e05ce9408c9d1c2230969.de-de.xlf"
$de.Range("F3").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c3cc36ee36b707a5742c80145a4eade4972fa2c9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/gt/a634b5f3-a252-4698-b996-c9ad1c439b66.1c5614d47b2feca54d8e05ce9408c9d1c2230969.de-de.xlf", "", "", "a634b5f3-a252-4698-b996-c9ad1c439b66.1c5614d47b2feca54d8e05ce9408c9d1c2230969.de-de.xlf")

$de.Range("G3").Value = "2016-02-18 09:51:25"
